# Italy commercial model update: correct "S/LFM/CDL" taxonomy code to "S/LFM+CDL"
# wherever it appears in the mapping table, and restore the active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the LFM/CDL -> LFM+CDL typo across all cells on the sheet (Offices, Trade,
# Hotels columns each contain one occurrence inside their multi-line mapping text).
$ws.Cells.Replace("S/LFM/CDL/", "S/LFM+CDL/")

# Update the saved selection/active cell to C2.
$ws.Range("C2").Select()
